$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(2,7,"155.8700226666667"),
    @(2,8,"467.610068"),
    @(2,9,"0.4627663557222626"),
    @(2,10,"0.4864916976605717"),
    @(2,13,"121.928739"),
    @(2,14,"365.786217"),
    @(2,15,"0.2282232151508951"),
    @(2,16,"0.2419720431319445"),
    @(2,17,"19005.03531164808"),
    @(2,18,"171045.3178048327"),
    @(2,19,"0.1056140255665976"),
    @(2,20,"0.1177173900496568"),
    @(3,7,"155.8700226666667"),
    @(3,8,"467.610068"),
    @(3,9,"0.4627663557222626"),
    @(3,10,"0.4864916976605717"),
    @(3,15,"0.2768624053389947"),
    @(3,16,"0.2935413991166814"),
    @(3,17,"23055.40996982481"),
    @(3,18,"207498.6897284233"),
    @(3,19,"0.1281226063552265"),
    @(3,20,"0.1428054535899338"),
    @(4,7,"155.8700226666667"),
    @(4,8,"467.610068"),
    @(4,9,"0.4627663557222626"),
    @(4,10,"0.4864916976605717"),
    @(4,13,"83.50496933333334"),
    @(4,14,"250.514908"),
    @(4,15,"0.1563025480180701"),
    @(4,16,"0.1657186665504434"),
    @(4,17,"13015.92146276597"),
    @(4,18,"117143.2931648937"),
    @(4,19,"0.07233156053642624"),
    @(4,20,"0.08062075542417138"),
    @(5,7,"155.8700226666667"),
    @(5,8,"467.610068"),
    @(5,9,"0.4627663557222626"),
    @(5,10,"0.4864916976605717"),
    @(5,13,"91.06846250000001"),
    @(5,14,"182.136925"),
    @(5,15,"0.1704597085236707"),
    @(5,16,"0.1204857969594293"),
    @(5,17,"14194.84331409348"),
    @(5,18,"85169.05988456091"),
    @(5,19,"0.07888301811097821"),
    @(5,20,"0.05861533990677972"),
    @(6,7,"155.8700226666667"),
    @(6,8,"467.610068"),
    @(6,9,"0.4627663557222626"),
    @(6,10,"0.4864916976605717"),
    @(6,13,"89.83562999999999"),
    @(6,14,"269.50689"),
    @(6,15,"0.1681521229683693"),
    @(6,16,"0.1782820942415013"),
    @(6,17,"14002.68168437428"),
    @(6,18,"126024.1351593685"),
    @(6,19,"0.07781514515303403"),
    @(6,20,"0.08673275869002998"),
    @(7,9,"0.3897411505765819"),
    @(7,10,"0.4097225989911443"),
    @(7,13,"121.928739"),
    @(7,14,"365.786217"),
    @(7,15,"0.2282232151508951"),
    @(7,16,"0.2419720431319445"),
    @(7,17,"16006.01305068893"),
    @(7,18,"144054.1174562004"),
    @(7,19,"0.08894797846119666"),
    @(7,20,"0.09914141439521758"),
    @(8,9,"0.3897411505765819"),
    @(8,10,"0.4097225989911443"),
    @(8,15,"0.2768624053389947"),
    @(8,16,"0.2935413991166814"),
    @(8,19,"0.1079046724082198"),
    @(8,20,"0.1202705449575835"),
    @(9,9,"0.3897411505765819"),
    @(9,10,"0.4097225989911443"),
    @(9,13,"83.50496933333334"),
    @(9,14,"250.514908"),
    @(9,15,"0.1563025480180701"),
    @(9,16,"0.1657186665504434"),
    @(9,17,"10961.99009280915"),
    @(9,18,"98657.91083528234"),
    @(9,19,"0.06091753490261407"),
    @(9,20,"0.06789868276039446"),
    @(10,9,"0.3897411505765819"),
    @(10,10,"0.4097225989911443"),
    @(10,13,"91.06846250000001"),
    @(10,14,"182.136925"),
    @(10,15,"0.1704597085236707"),
    @(10,16,"0.1204857969594293"),
    @(10,17,"11954.87635840453"),
    @(10,18,"71729.25815042715"),
    @(10,19,"0.06643516292696422"),
    @(10,20,"0.0493657538717367"),
    @(11,9,"0.3897411505765819"),
    @(11,10,"0.4097225989911443"),
    @(11,13,"89.83562999999999"),
    @(11,14,"269.50689"),
    @(11,15,"0.1681521229683693"),
    @(11,16,"0.1782820942415013"),
    @(11,17,"11793.03811381878"),
    @(11,18,"106137.343024369"),
    @(11,19,"0.06553580187758713"),
    @(11,20,"0.07304620300621202"),
    @(12,7,"0.2461213333333333"),
    @(12,8,"0.738364"),
    @(12,9,"0.0007307156985262189"),
    @(12,10,"0.0007681784042585035"),
    @(12,13,"121.928739"),
    @(12,14,"365.786217"),
    @(12,15,"0.2282232151508951"),
    @(12,16,"0.2419720431319445"),
    @(12,17,"30.009263814332"),
    @(12,18,"270.083374328988"),
    @(12,19,"0.0001667662860788858"),
    @(12,20,"0.0001858776979682669"),
    @(13,7,"0.2461213333333333"),
    @(13,8,"0.738364"),
    @(13,9,"0.0007307156985262189"),
    @(13,10,"0.0007681784042585035"),
    @(13,15,"0.2768624053389947"),
    @(13,16,"0.2935413991166814"),
    @(13,17,"36.40487211870666"),
    @(13,18,"327.64384906836"),
    @(13,19,"0.0002023077059129327"),
    @(13,20,"0.0002254921635572608"),
    @(14,7,"0.2461213333333333"),
    @(14,8,"0.738364"),
    @(14,9,"0.0007307156985262189"),
    @(14,10,"0.0007681784042585035"),
    @(14,13,"83.50496933333334"),
    @(14,14,"250.514908"),
    @(14,15,"0.1563025480180701"),
    @(14,16,"0.1657186665504434"),
    @(14,17,"20.55235439227911"),
    @(14,18,"184.971189530512"),
    @(14,19,"0.0001142127255564519"),
    @(14,20,"0.0001273015008265666"),
    @(15,7,"0.2461213333333333"),
    @(15,8,"0.738364"),
    @(15,9,"0.0007307156985262189"),
    @(15,10,"0.0007681784042585035"),
    @(15,13,"91.06846250000001"),
    @(15,14,"182.136925"),
    @(15,15,"0.1704597085236707"),
    @(15,16,"0.1204857969594293"),
    @(15,17,"22.41389141511667"),
    @(15,18,"134.4833484907"),
    @(15,19,"0.0001245575849844497"),
    @(15,20,"9.255458724410848E-05"),
    @(16,7,"0.2461213333333333"),
    @(16,8,"0.738364"),
    @(16,9,"0.0007307156985262189"),
    @(16,10,"0.0007681784042585035"),
    @(16,13,"89.83562999999999"),
    @(16,14,"269.50689"),
    @(16,15,"0.1681521229683693"),
    @(16,16,"0.1782820942415013"),
    @(16,17,"22.11046503644"),
    @(16,18,"198.99418532796"),
    @(16,19,"0.0001228713959934986"),
    @(16,20,"0.0001369524546623006"),
    @(17,7,"49.2786865"),
    @(17,8,"98.557373"),
    @(17,9,"0.1463047080910041"),
    @(17,10,"0.1025370217386683"),
    @(17,13,"121.928739"),
    @(17,14,"365.786217"),
    @(17,15,"0.2282232151508951"),
    @(17,16,"0.2419720431319445"),
    @(17,17,"6008.488104521323"),
    @(17,18,"36050.92862712793"),
    @(17,19,"0.03339013087224213"),
    @(17,20,"0.02481109264677019"),
    @(18,7,"49.2786865"),
    @(18,8,"98.557373"),
    @(18,9,"0.1463047080910041"),
    @(18,10,"0.1025370217386683"),
    @(18,15,"0.2768624053389947"),
    @(18,16,"0.2935413991166814"),
    @(18,17,"7289.023896927544"),
    @(18,18,"43734.14338156526"),
    @(18,19,"0.04050627339449488"),
    @(18,20,"0.03009886082242628"),
    @(19,7,"49.2786865"),
    @(19,8,"98.557373"),
    @(19,9,"0.1463047080910041"),
    @(19,10,"0.1025370217386683"),
    @(19,13,"83.50496933333334"),
    @(19,14,"250.514908"),
    @(19,15,"0.1563025480180701"),
    @(19,16,"0.1657186665504434"),
    @(19,17,"4115.015204969447"),
    @(19,18,"24690.09122981668"),
    @(19,19,"0.02286779866166389"),
    @(19,20,"0.01699229851458594"),
    @(20,7,"49.2786865"),
    @(20,8,"98.557373"),
    @(20,9,"0.1463047080910041"),
    @(20,10,"0.1025370217386683"),
    @(20,13,"91.06846250000001"),
    @(20,14,"182.136925"),
    @(20,15,"0.1704597085236707"),
    @(20,16,"0.1204857969594293"),
    @(20,17,"4487.734213574507"),
    @(20,18,"17950.93685429803"),
    @(20,19,"0.02493905789683329"),
    @(20,20,"0.01235425478202978"),
    @(21,7,"49.2786865"),
    @(21,8,"98.557373"),
    @(21,9,"0.1463047080910041"),
    @(21,10,"0.1025370217386683"),
    @(21,13,"89.83562999999999"),
    @(21,14,"269.50689"),
    @(21,15,"0.1681521229683693"),
    @(21,16,"0.1782820942415013"),
    @(21,17,"4426.981847299995"),
    @(21,18,"26561.89108379997"),
    @(21,19,"0.02460144726576989"),
    @(21,20,"0.01828051497285614"),
    @(22,7,"0.1539513333333333"),
    @(22,8,"0.461854"),
    @(22,9,"0.0004570699116250634"),
    @(22,10,"0.0004805032053572586"),
    @(22,13,"121.928739"),
    @(22,14,"365.786217"),
    @(22,15,"0.2282232151508951"),
    @(22,16,"0.2419720431319445"),
    @(22,17,"18.771091940702"),
    @(22,18,"168.939827466318"),
    @(22,19,"0.0001043139647798074"),
    @(22,20,"0.0001162683423317442"),
    @(23,7,"0.1539513333333333"),
    @(23,8,"0.461854"),
    @(23,9,"0.0004570699116250634"),
    @(23,10,"0.0004805032053572586"),
    @(23,15,"0.2768624053389947"),
    @(23,16,"0.2935413991166814"),
    @(23,17,"22.77160832260666"),
    @(23,18,"204.94447490346"),
    @(23,19,"0.0001265454751405968"),
    @(23,20,"0.0001410475831806198"),
    @(24,7,"0.1539513333333333"),
    @(24,8,"0.461854"),
    @(24,9,"0.0004570699116250634"),
    @(24,10,"0.0004805032053572586"),
    @(24,13,"83.50496933333334"),
    @(24,14,"250.514908"),
    @(24,15,"0.1563025480180701"),
    @(24,16,"0.1657186665504434"),
    @(24,17,"12.85570136882578"),
    @(24,18,"115.701312319432"),
    @(24,19,"7.14411918093915E-05"),
    @(24,20,"7.962835046501874E-05"),
    @(25,7,"0.1539513333333333"),
    @(25,8,"0.461854"),
    @(25,9,"0.0004570699116250634"),
    @(25,10,"0.0004805032053572586"),
    @(25,13,"91.06846250000001"),
    @(25,14,"182.136925"),
    @(25,15,"0.1704597085236707"),
    @(25,16,"0.1204857969594293"),
    @(25,17,"14.02011122649167"),
    @(25,18,"84.12066735895"),
    @(25,19,"7.791200391054825E-05"),
    @(25,20,"5.789381163902963E-05"),
    @(26,7,"0.1539513333333333"),
    @(26,8,"0.461854"),
    @(26,9,"0.0004570699116250634"),
    @(26,10,"0.0004805032053572586"),
    @(26,13,"89.83562999999999"),
    @(26,14,"269.50689"),
    @(26,15,"0.1681521229683693"),
    @(26,16,"0.1782820942415013"),
    @(26,17,"13.830315019339997"),
    @(26,18,"124.47283517406"),
    @(26,19,"7.685727598471934E-05"),
    @(26,20,"8.566511774084622E-05")
)

foreach ($item in $changes) {
    $r = $item[0]
    $c = $item[1]
    $v = [double]$item[2]
    $ws.Cells.Item($r, $c).Value = $v
}

Write-Output ("Applied " + $changes.Count + " cell updates")
